# contactDetail upload - final corrections pass.
# Fixes a couple of typo'd contact fields and tidies the trailing commas
# left over in the "Hobbies" column, then nudges the header row styling
# (row height / font colour) to match the final reviewed look.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2 (first contact) --------------------------------------------
$ws.Range("B2").Value = "rohan"               # FirstName: tinu -> rohan
$ws.Range("J2").Value = "tintu@gmail.com"     # Email: fix typo tinu@ -> tintu@
$ws.Range("L2").Value = "Reading ,Drawing"    # Hobbies: drop trailing comma

# --- Row 3 (second contact) --------------------------------------------
$ws.Range("B3").Value = "mini"                # FirstName: Maya -> mini
$ws.Range("H3").Value = "abcd"                # Street: dfbdf -> abcd
$ws.Range("L3").Value = "Reading ,Writing"    # Hobbies: drop trailing comma

# --- Formatting polish ---------------------------------------------------
# Pincode / Phone columns: make the font colour explicit black instead of
# the inherited theme colour.
$ws.Range("I2").Font.Color = 0
$ws.Range("K2").Font.Color = 0
$ws.Range("I3").Font.Color = 0
$ws.Range("K3").Font.Color = 0

# Slightly taller rows for the header + two data rows.
$ws.Range("A1:L3").Rows.RowHeight = 19.5
